# repull data, push all data, mean calculation
# Updates specific values in column F (dSF) to reflect repulled/pushed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = -2
    4  = -5
    9  = 4
    13 = -2
    14 = 3
    28 = 1
    36 = 0
    39 = -2
    48 = -3
    56 = -2
    60 = 1
    64 = -3
    67 = -1
    69 = -2
    70 = -8
    72 = 3
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
